$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Instance data should be stored in the scene, not with the model.
# Remove the "Instance Buffer Desc" row and the "Instance buffer data" row
# (delete bottom row first so the earlier row index stays valid).
$ws.Rows(26).Delete()
$ws.Rows(22).Delete()

# File format version bumped (no instance-buffer flag bit anymore).
$ws.Range("C1").Value = "Version 2"

# The old HAS_INSTANCE_BUFFER (0x0002) flag slot is reused for
# HAS_STRING_TABLE, which used to be 0x0008 (now 0x0001 since the two
# instance-buffer flag bits were removed).
$ws.Range("F5").Value = "HAS_STRING_TABLE"
$ws.Range("G5").Value = "0x0001"

# HAS_INSTANCE_BUFFER_DATA (0x0004) flag row is gone entirely.
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()

# Old HAS_STRING_TABLE (0x0008) row is now blank (its content moved to F5/G5 above).
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()

# Shrink the conditional-formatting range to match the two fewer rows.
$cf = $ws.Cells.FormatConditions
for ($i = 1; $i -le $cf.Count; $i++) {
    $cf.Item($i).ModifyAppliesToRange($ws.Range("C5:C119"))
}

# Match the author's final cursor position.
$ws.Range("A25:XFD25").Select()
